$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 169.16667
$ws.Range("I98").Value = 193
$ws.Range("K98").Value = 193
$ws.Range("M98").Value = 1305
$ws.Range("H122").Value = 169.16667
$ws.Range("I122").Value = 193
$ws.Range("K122").Value = 579
$ws.Range("M122").Value = 1871

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 954.8823
$ws.Range("J2").Value = 611.6667
$ws.Range("L2").Value = 611.6667
$ws.Range("N2").Value = -837.6667
$ws.Range("H8").Value = 8209.6
$ws.Range("J8").Value = 9166
$ws.Range("L8").Value = 9166
$ws.Range("N8").Value = -9454
$ws.Range("H11").Value = 4749.5
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 4749.5
$ws.Range("K11").Value = 0
$ws.Range("L11").ClearContents()
$ws.Range("M11").Value = 4749.5
$ws.Range("N11").Value = -5037.5
$ws.Range("H13").Value = 4350
$ws.Range("I13").Value = 1950
$ws.Range("J13").Value = 6750
$ws.Range("K13").Value = 1950
$ws.Range("L13").Value = 6750
$ws.Range("M13").Value = -1806
$ws.Range("N13").Value = -7038
$ws.Range("H14").Value = 3587.25
$ws.Range("I14").Value = 3175
$ws.Range("K14").Value = 3175
$ws.Range("M14").Value = -3000
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").ClearContents()
$ws.Range("N15").Value = 0
$ws.Range("H19").Value = 5251.1665
$ws.Range("I19").Value = 6627
$ws.Range("K19").Value = 6627
$ws.Range("M19").Value = -6398
$ws.Range("H21").Value = 7407
$ws.Range("I21").Value = 12242.5
$ws.Range("J21").Value = 4183.3335
$ws.Range("K21").Value = 12242.5
$ws.Range("L21").Value = 4183.3335
$ws.Range("M21").Value = -11868.5
$ws.Range("N21").Value = -4931.3335
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").ClearContents()
$ws.Range("N23").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").ClearContents()
$ws.Range("M25").ClearContents()
$ws.Range("N25").Value = 0
$ws.Range("H30").Value = 4499.6665
$ws.Range("I30").Value = 1000
$ws.Range("J30").Value = 6249.5
$ws.Range("K30").Value = 1000
$ws.Range("L30").Value = 6249.5
$ws.Range("M30").Value = -850
$ws.Range("N30").Value = -6549.5
$ws.Range("H36").Value = 2026
$ws.Range("I36").Value = 2026
$ws.Range("K36").Value = 2026
$ws.Range("M36").Value = -1680
$ws.Range("H38").Value = 1699
$ws.Range("I38").Value = 998.75
$ws.Range("J38").Value = 4500
$ws.Range("K38").Value = 998.75
$ws.Range("L38").Value = 4500
$ws.Range("M38").Value = -531.75
$ws.Range("N38").Value = -5434
$ws.Range("H40").Value = 5000
$ws.Range("J40").Value = 5000
$ws.Range("L40").Value = 5000
$ws.Range("N40").Value = -5352
$ws.Range("H97").Value = 530.35297
$ws.Range("I97").Value = 577.0909
$ws.Range("J97").Value = 444.66666
$ws.Range("K97").Value = 577.0909
$ws.Range("L97").Value = 444.66666
$ws.Range("M97").Value = -81.09090000000003
$ws.Range("N97").Value = -1436.66666
$ws.Range("H116").Value = 954.8823
$ws.Range("J116").Value = 611.6667
$ws.Range("L116").Value = 611.6667
$ws.Range("N116").Value = -5199.6667

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 954.8823
$ws.Range("J3").Value = 611.6667
$ws.Range("L3").Value = 611.6667
$ws.Range("N3").Value = -839.6667
$ws.Range("H22").Value = 9001
$ws.Range("I22").Value = 9001
$ws.Range("K22").Value = 9001
$ws.Range("M22").Value = -8828
$ws.Range("H25").Value = 1866
$ws.Range("I25").Value = 2304.6667
$ws.Range("J25").Value = 550
$ws.Range("K25").Value = 2304.6667
$ws.Range("L25").Value = 550
$ws.Range("M25").Value = -2069.6667
$ws.Range("N25").Value = -1020
$ws.Range("H29").Value = 558.3333
$ws.Range("I29").Value = 558.3333
$ws.Range("K29").Value = 558.3333
$ws.Range("M29").Value = -269.3333
$ws.Range("H30").Value = 300
$ws.Range("J30").Value = 300
$ws.Range("L30").Value = 300
$ws.Range("N30").Value = -550
$ws.Range("H36").Value = 2024.3334
$ws.Range("I36").Value = 2024.3334
$ws.Range("K36").Value = 2024.3334
$ws.Range("M36").Value = -1490.3334
$ws.Range("H37").Value = 1233.8334
$ws.Range("I37").Value = 1233.8334
$ws.Range("K37").Value = 1233.8334
$ws.Range("M37").Value = -1096.8334
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").ClearContents()
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = 0

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 455.9
$ws.Range("J5").Value = 500
$ws.Range("L5").Value = 500
$ws.Range("N5").Value = -724
$ws.Range("H31").Value = 8377.789000000001
$ws.Range("I31").Value = 4538
$ws.Range("K31").Value = 4538
$ws.Range("M31").Value = -4243
$ws.Range("H34").Value = 8377.789000000001
$ws.Range("I34").Value = 4538
$ws.Range("K34").Value = 4538
$ws.Range("M34").Value = -4336
$ws.Range("H55").Value = 6036
$ws.Range("I55").Value = 6036
$ws.Range("K55").Value = 6036
$ws.Range("M55").Value = -5721
$ws.Range("H56").Value = 292.5
$ws.Range("I56").Value = 292.5
$ws.Range("K56").Value = 292.5
$ws.Range("M56").Value = 552.5
$ws.Range("H94").Value = 5052.3335
$ws.Range("J94").Value = 7303
$ws.Range("L94").Value = 7303
$ws.Range("N94").Value = -8205
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").ClearContents()
$ws.Range("N98").Value = 0
$ws.Range("H99").Value = 1921.4
$ws.Range("I99").Value = 1864.2
$ws.Range("K99").Value = 1864.2
$ws.Range("M99").Value = -366.2
$ws.Range("H126").Value = 1921.4
$ws.Range("I126").Value = 1864.2
$ws.Range("K126").Value = 5592.6
$ws.Range("M126").Value = -3122.6

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 895.6
$ws.Range("I5").Value = 857.6
$ws.Range("J5").Value = 933.6
$ws.Range("K5").Value = 2572.8
$ws.Range("L5").Value = 2800.8
$ws.Range("M5").Value = -2460.8
$ws.Range("N5").Value = -3024.8
$ws.Range("H41").Value = 4000
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 4000
$ws.Range("K41").Value = 0
$ws.Range("L41").ClearContents()
$ws.Range("M41").Value = 12000
$ws.Range("N41").Value = -12676
$ws.Range("H135").Value = 895.6
$ws.Range("I135").Value = 857.6
$ws.Range("J135").Value = 933.6
$ws.Range("K135").Value = 7718.400000000001
$ws.Range("L135").Value = 8402.4
$ws.Range("M135").Value = -5183.400000000001
$ws.Range("N135").Value = -13472.4

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 372
$ws.Range("I97").Value = 319
$ws.Range("J97").Value = 467.4
$ws.Range("K97").Value = 319
$ws.Range("L97").Value = 467.4
$ws.Range("M97").Value = 177
$ws.Range("N97").Value = -1459.4

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2100
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 2100
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
